# Adding parameterization, creating third test
$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("AddCustomerTest")
$wsOpen = $wb.Worksheets.Item("OpenAccountTest")

# --- OpenAccountTest: add new "alerttext" column first, so the new shared
#     string "Account created successfully" is registered before the
#     AddCustomerTest parameterized rows' new strings. ---
$wsOpen.Range("C2").Value = "Account created successfully"
$wsOpen.Range("C1").Value = "alerttext"
$wsOpen.Columns.Item(3).ColumnWidth = 26.140625

# --- AddCustomerTest: add 3 more parameterized rows ---
$wsAdd.Range("A3").Value = "Petya"
$wsAdd.Range("B3").Value = "Petrov"
$wsAdd.Range("C3").Value = "asdasd"
$wsAdd.Range("D3").Value = "Customer added successfully"

$wsAdd.Range("A4").Value = "Sidor"
$wsAdd.Range("B4").Value = "Sidorov"
$wsAdd.Range("C4").Value = "sdsagg4"
$wsAdd.Range("D4").Value = "Customer added successfully"

$wsAdd.Range("A5").Value = "Kirill"
$wsAdd.Range("B5").Value = "Kirillov"
$wsAdd.Range("C5").Value = "sdfgre34"
$wsAdd.Range("D5").Value = "Customer added successfully"

$wsAdd.Range("D4:D5").Select() | Out-Null

$wsOpen.Range("C4").Select() | Out-Null

# --- Switch active tab to AddCustomerTest ---
$wsAdd.Select() | Out-Null
